# Updates the crypto price/volume table (and two pairs of swapped rows)
# to match the latest scraped values from coinranking.com.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.223.13"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "3.404.38"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.00"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.18"
$ws.Range("E6").Value = "  +1.78%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.404.47"
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.72"
$ws.Range("E10").Value = "  +3.41%  "
$ws.Range("E11").Value = "  -1.25%  "
$ws.Range("E12").Value = "  -1.87%  "
$ws.Range("D13").Value = "3.984.22"
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.64"
$ws.Range("E15").Value = "  +2.03%  "
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").Value = "3.403.09"
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("D18").Value = "61.305.79"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("E19").Value = "  +1.34%  "
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.33"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "378.64"
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("B23").Value = "WrappedeETH"
$ws.Range("C23").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D23").Value = "3.536.51"
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.553"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("E28").Value = "  +10.05%  "
$ws.Range("E29").Value = "  -7.51%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.44"
$ws.Range("E30").Value = "  -1.17%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.17"
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.56"
$ws.Range("E36").Value = "  +1.77%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.91"
$ws.Range("E37").Value = "  +1.33%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.13"
$ws.Range("E38").Value = "  -2.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "166.25"
$ws.Range("E39").Value = "  +0.54%  "
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.12"
$ws.Range("E41").Value = "  +6.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.75"
$ws.Range("E42").Value = "  +1.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  +0.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.98"
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.17"
$ws.Range("E47").Value = "  -2.68%  "
$ws.Range("D48").Value = "2.539.45"
$ws.Range("E48").Value = "  +7.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.89"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("E51").Value = "  -0.10%  "

$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
